# Update Betfair Back/Lay odds values for 2026-01-03 on Sheet1.
# Applies the 174 individual cell value updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 4.1
$ws.Range("K2").Value = 4.2
$ws.Range("L2").Value = 1.37
$ws.Range("O2").Value = 1.27
$ws.Range("P2").Value = 2.04
$ws.Range("Q2").Value = 1.81
$ws.Range("R2").Value = 1.4
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 1.79
$ws.Range("W2").Value = 2.28
$ws.Range("X2").Value = 16
$ws.Range("Y2").Value = 19.5
$ws.Range("Z2").Value = 40
$ws.Range("AA2").Value = 130
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 20
$ws.Range("AE2").Value = 70
$ws.Range("AF2").Value = 11
$ws.Range("AG2").Value = 10
$ws.Range("AH2").Value = 19.5
$ws.Range("AI2").Value = 70
$ws.Range("AJ2").Value = 18.5
$ws.Range("AK2").Value = 18
$ws.Range("AL2").Value = 34
$ws.Range("AM2").Value = 110
$ws.Range("AN2").Value = 12.5
$ws.Range("AO2").Value = 75
$ws.Range("L3").Value = 1.31
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 2.08
$ws.Range("Q3").Value = 1.84
$ws.Range("R3").Value = 1.43
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 1.69
$ws.Range("U3").Value = 2.26
$ws.Range("AN3").Value = 16.5
$ws.Range("AO3").Value = 34
$ws.Range("N6").Value = 4.9
$ws.Range("O6").Value = 1.21
$ws.Range("R6").Value = 1.56
$ws.Range("T6").Value = 1.6
$ws.Range("U6").Value = 2.4
$ws.Range("X6").Value = 1000
$ws.Range("AH6").Value = 16.5
$ws.Range("AK6").Value = 1000
$ws.Range("F9").Value = 1.85
$ws.Range("P9").Value = 1.77
$ws.Range("Q9").Value = 2.06
$ws.Range("F10").Value = 2.68
$ws.Range("H10").Value = 2.46
$ws.Range("J10").Value = 3.3
$ws.Range("Q10").Value = 1.67
$ws.Range("Q11").Value = 2.72
$ws.Range("F14").Value = 1.96
$ws.Range("G14").Value = 2.62
$ws.Range("H14").Value = 3.2
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = 5.9
$ws.Range("J16").Value = 3.45
$ws.Range("P16").Value = 1.87
$ws.Range("N17").Value = 3.3
$ws.Range("T17").Value = 1.84
$ws.Range("AB17").Value = 8.800000000000001
$ws.Range("F18").Value = 1.92
$ws.Range("G18").Value = 2.08
$ws.Range("H18").Value = 3.8
$ws.Range("I18").Value = 4.3
$ws.Range("AC18").Value = 11.5
$ws.Range("Q19").Value = 1.9
$ws.Range("AC19").Value = 9.199999999999999
$ws.Range("AJ19").Value = 1000
$ws.Range("S20").Value = 2.9
$ws.Range("Z20").Value = 65
$ws.Range("H21").Value = 2.88
$ws.Range("I21").Value = 2.96
$ws.Range("AA21").Value = 60
$ws.Range("AJ21").Value = 40
$ws.Range("AK21").Value = 28
$ws.Range("AN21").Value = 24
$ws.Range("AO21").Value = 36
$ws.Range("H22").Value = 7.2
$ws.Range("P22").Value = 2.16
$ws.Range("Q22").Value = 1.69
$ws.Range("G23").Value = 2.64
$ws.Range("H23").Value = 2.64
$ws.Range("H24").Value = 8.800000000000001
$ws.Range("J24").Value = 5.1
$ws.Range("N24").Value = 4.5
$ws.Range("P24").Value = 2.22
$ws.Range("F25").Value = 1.7
$ws.Range("G25").Value = 1.88
$ws.Range("I25").Value = 5.6
$ws.Range("G26").Value = 1.87
$ws.Range("H26").Value = 3.95
$ws.Range("F27").Value = 2.38
$ws.Range("G28").Value = 2.9
$ws.Range("I28").Value = 3.95
$ws.Range("J28").Value = 2.94
$ws.Range("F29").Value = 2.74
$ws.Range("I29").Value = 2.72
$ws.Range("Q31").Value = 1.66
$ws.Range("P32").Value = 1.91
$ws.Range("Q32").Value = 1.87
$ws.Range("F33").Value = 1.51
$ws.Range("G33").Value = 1.67
$ws.Range("H33").Value = 5.6
$ws.Range("I33").Value = 7.8
$ws.Range("K33").Value = 6.2
$ws.Range("P33").Value = 2.54
$ws.Range("Q33").Value = 1.46
$ws.Range("S34").Value = 5.2
$ws.Range("U34").Value = 1.86
$ws.Range("Z34").Value = 16.5
$ws.Range("AD34").Value = 13.5
$ws.Range("F35").Value = 1.75
$ws.Range("G35").Value = 1.86
$ws.Range("H35").Value = 4.8
$ws.Range("I35").Value = 5.6
$ws.Range("K35").Value = 4.4
$ws.Range("Q35").Value = 1.66
$ws.Range("F36").Value = 1.55
$ws.Range("G36").Value = 1.66
$ws.Range("Q36").Value = 1.79
$ws.Range("F37").Value = 2.78
$ws.Range("G37").Value = 36
$ws.Range("H37").Value = 1.31
$ws.Range("I37").Value = 1.56
$ws.Range("J37").Value = 4.8
$ws.Range("Q37").Value = 1.62
$ws.Range("Q38").Value = 2.02
$ws.Range("H39").Value = 3.8
$ws.Range("P39").Value = 2.48
$ws.Range("Q39").Value = 1.66
$ws.Range("R39").Value = 1.59
$ws.Range("S39").Value = 2.54
$ws.Range("T39").Value = 1.59
$ws.Range("Z39").Value = 32
$ws.Range("F41").Value = 1.61
$ws.Range("I41").Value = 11
$ws.Range("J41").Value = 3.5
$ws.Range("P41").Value = 1.87
$ws.Range("Q41").Value = 1.7
$ws.Range("F43").Value = 7
$ws.Range("G43").Value = 14
$ws.Range("K43").Value = 6.4
$ws.Range("Q43").Value = 1.53
$ws.Range("AE44").Value = 1000
$ws.Range("N46").Value = 4.9
$ws.Range("U46").Value = 2.12
$ws.Range("Z46").Value = 9.4
$ws.Range("F47").Value = 1.82
$ws.Range("I47").Value = 5.6
$ws.Range("P47").Value = 1.96
$ws.Range("Q47").Value = 1.82
$ws.Range("O48").Value = 1.26
$ws.Range("R48").Value = 1.46
$ws.Range("S48").Value = 3
$ws.Range("H49").Value = 3.3
$ws.Range("I49").Value = 3.6
$ws.Range("Q49").Value = 1.66
$ws.Range("P50").Value = 2.1
$ws.Range("Q50").Value = 1.83
$ws.Range("AF50").Value = 22
$ws.Range("I51").Value = 38
$ws.Range("P51").Value = 2.64
$ws.Range("Q51").Value = 1.44
$ws.Range("G54").Value = 6.4
$ws.Range("F55").Value = 2.1
$ws.Range("I55").Value = 4.5
$ws.Range("Q55").Value = 2.48
$ws.Range("S56").Value = 2.8
$ws.Range("AE56").Value = 980
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 2.24
